# Apply cryptocurrency price/volume updates scraped on Tue Apr 23 12:14:22 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.019.46"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "3.179.49"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'605.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.69%  "
$ws.Range("D6").Value = "'154.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.180.09"
$ws.Range("E8").Value = "  -0.55%  "
$ws.Range("E9").Value = "  +1.79%  "
$ws.Range("E10").Value = "  -1.12%  "
$ws.Range("D11").Value = "'5.62"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.54%  "
$ws.Range("E12").Value = "  +0.56%  "
$ws.Range("D13").Value = "'0.0000265"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.20%  "
$ws.Range("D14").Value = "'38.22"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.20%  "
$ws.Range("D15").Value = "3.701.49"
$ws.Range("E15").Value = "  -0.55%  "
$ws.Range("D16").Value = "66.109.58"
$ws.Range("E16").Value = "  +0.25%  "
$ws.Range("D17").Value = "'7.38"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.71%  "
$ws.Range("D18").Value = "3.181.92"
$ws.Range("E18").Value = "  -0.69%  "
$ws.Range("E19").Value = "  +0.91%  "
$ws.Range("D20").Value = "'506.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.78%  "
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("D22").Value = "'0.729"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.64%  "
$ws.Range("D23").Value = "'7.97"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.19%  "
$ws.Range("D24").Value = "'14.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.53%  "
$ws.Range("D25").Value = "'84.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.75%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("D27").Value = "'3.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.23%  "
$ws.Range("D28").Value = "'9.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.10%  "
$ws.Range("D29").Value = "'2.38"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.98%  "
$ws.Range("B30").Value = "Stacks"
$ws.Range("C30").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D30").Value = "'3.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.99%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "'7.18"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.00%  "
$ws.Range("D32").Value = "'27.87"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.80%  "
$ws.Range("E33").Value = "  +0.14%  "
$ws.Range("E34").Value = "  -3.35%  "
$ws.Range("D35").Value = "'6.48"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.17%  "
$ws.Range("D36").Value = "'513.94"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.75%  "
$ws.Range("D37").Value = "'55.40"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.78%  "
$ws.Range("D38").Value = "'0.0882"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.49%  "
$ws.Range("D39").Value = "'0.0419"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.03%  "
$ws.Range("E40").Value = "  +5.06%  "
$ws.Range("D41").Value = "'8.74"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.79%  "
$ws.Range("D42").Value = "0.0₃0683"
$ws.Range("E42").Value = "  +5.59%  "
$ws.Range("D43").Value = "'2.86"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.05%  "
$ws.Range("D44").Value = "'0.297"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.06%  "
$ws.Range("D45").Value = "'2.46"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.16%  "
$ws.Range("D46").Value = "2.824.07"
$ws.Range("E46").Value = "  -3.70%  "
$ws.Range("D47").Value = "'27.90"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.74%  "
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("E49").Value = "  +2.34%  "
$ws.Range("E50").Value = "  +0.12%  "
$ws.Range("D51").Value = "'2.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.44%  "
